$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write order reproduces the shared-string table order seen in the target file.
$ws.Range("D2").Value = "TEST Adv Qty"
$ws.Range("A2:A5").Value = "Add"
$ws.Range("B2").Value = "Testing Advanced Quantities"
$ws.Range("C2").Value = "Testing Advanced Quantities"
$ws.Range("V3").Value = "y"
$ws.Range("AA2").Value = "10..100[10]"
$ws.Range("B3").Value = "Testing Any Quantity"
$ws.Range("C3").Value = "Testing Any Quantity"
$ws.Range("D3").Value = "TEST Any Qty"
$ws.Range("E2").Value = "Order only in quantities from 10 to 100, in steps of 10."
$ws.Range("E3").Value = "Order any quantity."
$ws.Range("B4").Value = "Testing Minimum Quantity"
$ws.Range("C4").Value = "Testing Minimum Quantity"
$ws.Range("D4").Value = "TEST Min Qty"
$ws.Range("E4").Value = "Order minimum of 250."
$ws.Range("B5").Value = "Testing Multiple Quantity"
$ws.Range("C5").Value = "Testing Multiple Quantity"
$ws.Range("D5").Value = "TEST Mult Qty"
$ws.Range("E5").Value = "Order in multiples of 5."

$ws.Range("AE2").Value = 1
$ws.Range("AH2").Value = 1
$ws.Range("AJ2").Value = 1
$ws.Range("AL2").Value = 1
$ws.Range("AN2").Value = 1
$ws.Range("AO2").Value = 0.5
$ws.Range("AP2").Value = 0

$ws.Range("AE3").Value = 1
$ws.Range("AH3").Value = 1
$ws.Range("AJ3").Value = 1
$ws.Range("AL3").Value = 1
$ws.Range("AN3").Value = 1
$ws.Range("AO3").Value = 0.5
$ws.Range("AP3").Value = 0

$ws.Range("AE4").Value = 1
$ws.Range("AH4").Value = 1
$ws.Range("AJ4").Value = 1
$ws.Range("AL4").Value = 1
$ws.Range("AN4").Value = 1
$ws.Range("AO4").Value = 0.5
$ws.Range("AP4").Value = 0

$ws.Range("X5").Value = 5
$ws.Range("Z5").Value = 5
$ws.Range("AE5").Value = 1
$ws.Range("AH5").Value = 1
$ws.Range("AJ5").Value = 1
$ws.Range("AL5").Value = 1
$ws.Range("AN5").Value = 1
$ws.Range("AO5").Value = 0.5
$ws.Range("AP5").Value = 0

$ws.Range("AE4:AP5").Select()

# Extend the existing "Text Length <= 50" rule from B1:B1048576 to also
# cover C2:C5 (the newly populated Display Name cells). The host's
# Validation.Add only honours the first area of a multi-area Range, so
# the two areas are applied with separate calls carrying identical
# settings (net effect: every target cell enforces the same rule).
$ws.Range("B1:B1048576").Validation.Delete()

$ws.Range("B1:B1048576").Validation.Add(6, 1, 8, 50)
$bValidation = $ws.Range("B1:B1048576").Validation
$bValidation.IgnoreBlank = $true
$bValidation.ShowInput = $false
$bValidation.ShowError = $true
$bValidation.ErrorTitle = "Text Length"
$bValidation.ErrorMessage = "Must be 50 or fewer characters."
$bValidation.InputTitle = "Test"
$bValidation.InputMessage = "Test message"

$ws.Range("C2:C5").Validation.Add(6, 1, 8, 50)
$cValidation = $ws.Range("C2:C5").Validation
$cValidation.IgnoreBlank = $true
$cValidation.ShowInput = $false
$cValidation.ShowError = $true
$cValidation.ErrorTitle = "Text Length"
$cValidation.ErrorMessage = "Must be 50 or fewer characters."
$cValidation.InputTitle = "Test"
$cValidation.InputMessage = "Test message"
